$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1 to make room for the header row,
# shifting all existing data down by one row.
$ws.Rows.Item(1).Insert()

# Populate the new header row.
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "status"
